{"js": "// The commit adds a \"\u041f\u0440\u043e\u0442\u043e\u043a\u043e\u043b\u0430 \u2116 ICNUM \u043e\u0442 ICDATE\" continuation to the\n// bulletin title paragraph (as 4 new runs, two of them tagged en-US),\n// marks the paragraph-mark itself as en-US, and relocates the single\n// \"_GoBack\" bookmark from its old spot (an empty paragraph near the\n// signature block) to the end of the title paragraph.\n\n// 1) Drop the \"_GoBack\" bookmark from wherever it currently lives -\n//    bookmark names are unique, so it has to move rather than be\n//    duplicated.\nconst existingGoBack = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\nexistingGoBack.load(\"isNullObject\");\nawait context.sync();\nif (!existingGoBack.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// 2) Rebuild the first paragraph (the \"\u0411\u042e\u041b\u041b\u0415\u0422\u0415\u041d\u042c \" title line) with the\n//    extra runs, the en-US language tag on the paragraph mark, and the\n//    bookmark re-created at the very end of the paragraph.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titleParagraph = paragraphs.items[0];\nconst titleRange = titleParagraph.getRange();\n\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">\n<w:body>\n<w:p w14:paraId=\"5F9BA1F8\" w14:textId=\"2DB6943F\" w:rsidR=\"00382498\" w:rsidRPr=\"000566D8\" w:rsidRDefault=\"00382498\" w:rsidP=\"00F073F4\">\n<w:pPr>\n<w:pStyle w:val=\"BodyTextIndent2\"/>\n<w:tabs><w:tab w:val=\"left\" w:pos=\"284\"/></w:tabs>\n<w:spacing w:after=\"0\" w:line=\"240\" w:lineRule=\"auto\"/>\n<w:ind w:left=\"0\"/>\n<w:jc w:val=\"center\"/>\n<w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:b/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/><w:lang w:val=\"en-US\"/></w:rPr>\n</w:pPr>\n<w:r w:rsidRPr=\"000566D8\">\n<w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:b/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr>\n<w:t xml:space=\"preserve\">\u0411\u042e\u041b\u041b\u0415\u0422\u0415\u041d\u042c </w:t>\n</w:r>\n<w:r>\n<w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:b/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr>\n<w:t xml:space=\"preserve\">\u041f\u0440\u043e\u0442\u043e\u043a\u043e\u043b\u0430 \u2116 </w:t>\n</w:r>\n<w:r>\n<w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:b/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/><w:lang w:val=\"en-US\"/></w:rPr>\n<w:t xml:space=\"preserve\">ICNUM </w:t>\n</w:r>\n<w:r>\n<w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:b/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr>\n<w:t xml:space=\"preserve\">\u043e\u0442 </w:t>\n</w:r>\n<w:r>\n<w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:b/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/><w:lang w:val=\"en-US\"/></w:rPr>\n<w:t>ICDATE</w:t>\n</w:r>\n<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n<w:bookmarkEnd w:id=\"0\"/>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n\ntitleRange.insertOoxml(ooxml, \"Replace\");\nawait context.sync();\n", "ps1": "# The commit adds a \"\u041f\u0440\u043e\u0442\u043e\u043a\u043e\u043b\u0430 \u2116 ICNUM \u043e\u0442 ICDATE\" continuation to the\n# bulletin title paragraph (as 4 new runs, two of them tagged en-US),\n# marks the paragraph-mark itself as en-US, and relocates the single\n# \"_GoBack\" bookmark from its old spot (an empty paragraph near the\n# signature block) to the end of the title paragraph.\n\n$d = $word.ActiveDocument\n\n# 1) Drop the \"_GoBack\" bookmark from wherever it currently lives -\n#    bookmark names are unique, so it has to move rather than be\n#    duplicated.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2) Rebuild the first paragraph (the \"\u0411\u042e\u041b\u041b\u0415\u0422\u0415\u041d\u042c \" title line) with the\n#    extra runs, the en-US language tag on the paragraph mark, and the\n#    bookmark re-created at the very end of the paragraph. Range.InsertXML\n#    replaces the target range's content, so target the whole paragraph\n#    (including its end-of-paragraph mark) to rewrite it in place.\n$titleParagraph = $d.Paragraphs(1)\n$titleRange = $titleParagraph.Range\n\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' + `\n'<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n'<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n'<pkg:xmlData>' + `\n'<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' + `\n'<w:body>' + `\n'<w:p w14:paraId=\"5F9BA1F8\" w14:textId=\"2DB6943F\" w:rsidR=\"00382498\" w:rsidRPr=\"000566D8\" w:rsidRDefault=\"00382498\" w:rsidP=\"00F073F4\">' + `\n'<w:pPr><w:pStyle w:val=\"BodyTextIndent2\"/><w:tabs><w:tab w:val=\"left\" w:pos=\"284\"/></w:tabs><w:spacing w:after=\"0\" w:line=\"240\" w:lineRule=\"auto\"/><w:ind w:left=\"0\"/><w:jc w:val=\"center\"/><w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:b/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' + `\n'<w:r w:rsidRPr=\"000566D8\"><w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:b/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr><w:t xml:space=\"preserve\">\u0411\u042e\u041b\u041b\u0415\u0422\u0415\u041d\u042c </w:t></w:r>' + `\n'<w:r><w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:b/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr><w:t xml:space=\"preserve\">\u041f\u0440\u043e\u0442\u043e\u043a\u043e\u043b\u0430 \u2116 </w:t></w:r>' + `\n'<w:r><w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:b/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\">ICNUM </w:t></w:r>' + `\n'<w:r><w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:b/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr><w:t xml:space=\"preserve\">\u043e\u0442 </w:t></w:r>' + `\n'<w:r><w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:b/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t>ICDATE</w:t></w:r>' + `\n'<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' + `\n'</w:p>' + `\n'</w:body></w:document>' + `\n'</pkg:xmlData></pkg:part></pkg:package>'\n\n$titleRange.InsertXML($xml)\n"}
